$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 93
$ws1.Range("F6").Value = 39
$ws1.Range("F7").Value = 597
$ws1.Range("F9").Value = 8808
$ws1.Range("F10").Value = 815
$ws1.Range("F11").Value = 332
$ws1.Range("F12").Value = 1147
$ws1.Range("F13").Value = 998
$ws1.Range("F14").Value = 115
$ws1.Range("F17").Value = 237
$ws1.Range("F18").Value = 276
$ws1.Range("F19").Value = 69
$ws1.Range("F21").Value = 1063

# Sheet "全部类型" (sheet4): update column F (想去人数) for the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 93
$ws4.Range("F7").Value = 39
$ws4.Range("F9").Value = 597
$ws4.Range("F11").Value = 8808
$ws4.Range("F12").Value = 815
$ws4.Range("F13").Value = 332
$ws4.Range("F14").Value = 1147
$ws4.Range("F15").Value = 998
$ws4.Range("F16").Value = 115
$ws4.Range("F19").Value = 237
$ws4.Range("F20").Value = 276
$ws4.Range("F21").Value = 69
$ws4.Range("F23").Value = 1063
